$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7454756666666666
$ws.Range("H2").Value = 2.236427
$ws.Range("I2").Value = 0.004237455828433692
$ws.Range("J2").Value = 0.004251944035061194
$ws.Range("M2").Value = 0.346153
$ws.Range("N2").Value = 1.038459
$ws.Range("O2").Value = 0.002587513800919532
$ws.Range("P2").Value = 0.002593123140481257
$ws.Range("Q2").Value = 0.2580486384436667
$ws.Range("R2").Value = 2.322437745993
$ws.Range("S2").Value = 0.00001096447543685909
$ws.Range("T2").Value = 0.00001102581446934843
$ws.Range("G3").Value = 0.7454756666666666
$ws.Range("H3").Value = 2.236427
$ws.Range("I3").Value = 0.004237455828433692
$ws.Range("J3").Value = 0.004251944035061194
$ws.Range("O3").Value = 0.0001556606107424992
$ws.Range("P3").Value = 0.000155998059463248
$ws.Range("Q3").Value = 0.01552378528266667
$ws.Range("R3").Value = 0.139714067544
$ws.Range("S3").Value = 0.0000006596049622483514
$ws.Range("T3").Value = 0.0000006632950184158789
$ws.Range("G4").Value = 0.7454756666666666
$ws.Range("H4").Value = 2.236427
$ws.Range("I4").Value = 0.004237455828433692
$ws.Range("J4").Value = 0.004251944035061194
$ws.Range("M4").Value = 84.40796133333333
$ws.Range("N4").Value = 253.223884
$ws.Range("O4").Value = 0.630954418587972
$ws.Range("P4").Value = 0.6323222325801418
$ws.Range("Q4").Value = 62.92408124694089
$ws.Range("R4").Value = 566.316731222468
$ws.Range("S4").Value = 0.002673641478521593
$ws.Range("T4").Value = 0.002688598745055711
$ws.Range("G5").Value = 0.7454756666666666
$ws.Range("H5").Value = 2.236427
$ws.Range("I5").Value = 0.004237455828433692
$ws.Range("J5").Value = 0.004251944035061194
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.868151
$ws.Range("N5").Value = 1.736302
$ws.Range("O5").Value = 0.006489479200764093
$ws.Range("P5").Value = 0.004335698275101748
$ws.Range("Q5").Value = 0.6471854454923334
$ws.Range("R5").Value = 3.883112672954
$ws.Range("S5").Value = 0.00002749888146277702
$ws.Range("T5").Value = 0.00001843514641864399
$ws.Range("G6").Value = 0.7454756666666666
$ws.Range("H6").Value = 2.236427
$ws.Range("I6").Value = 0.004237455828433692
$ws.Range("J6").Value = 0.004251944035061194
$ws.Range("M6").Value = 48.13513433333333
$ws.Range("N6").Value = 144.405403
$ws.Range("O6").Value = 0.3598129277996019
$ws.Range("P6").Value = 0.3605929479448119
$ws.Range("Q6").Value = 35.88357135723122
$ws.Range("R6").Value = 322.9521422150809
$ws.Range("S6").Value = 0.001524691388050214
$ws.Range("T6").Value = 0.001533221034099075
$ws.Range("G7").Value = 102.35201
$ws.Range("H7").Value = 307.05603
$ws.Range("I7").Value = 0.5817924591230612
$ws.Range("J7").Value = 0.583781654929077
$ws.Range("M7").Value = 0.346153
$ws.Range("N7").Value = 1.038459
$ws.Range("O7").Value = 0.002587513800919532
$ws.Range("P7").Value = 0.002593123140481257
$ws.Range("Q7").Value = 35.42945531753
$ws.Range("R7").Value = 318.86509785777
$ws.Range("S7").Value = 0.001505396017251834
$ws.Range("T7").Value = 0.001513817718385033
$ws.Range("G8").Value = 102.35201
$ws.Range("H8").Value = 307.05603
$ws.Range("I8").Value = 0.5817924591230612
$ws.Range("J8").Value = 0.583781654929077
$ws.Range("O8").Value = 0.0001556606107424992
$ws.Range("P8").Value = 0.000155998059463248
$ws.Range("Q8").Value = 2.13137825624
$ws.Range("R8").Value = 19.18240430616
$ws.Range("S8").Value = 0.00009056216951247623
$ws.Range("T8").Value = 0.00009106880531917948
$ws.Range("G9").Value = 102.35201
$ws.Range("H9").Value = 307.05603
$ws.Range("I9").Value = 0.5817924591230612
$ws.Range("J9").Value = 0.583781654929077
$ws.Range("M9").Value = 84.40796133333333
$ws.Range("N9").Value = 253.223884
$ws.Range("O9").Value = 0.630954418587972
$ws.Range("P9").Value = 0.6323222325801418
$ws.Range("Q9").Value = 8639.324502468946
$ws.Range("R9").Value = 77753.92052222051
$ws.Range("S9").Value = 0.3670845227848575
$ws.Range("T9").Value = 0.3691381193840839
$ws.Range("G10").Value = 102.35201
$ws.Range("H10").Value = 307.05603
$ws.Range("I10").Value = 0.5817924591230612
$ws.Range("J10").Value = 0.583781654929077
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.868151
$ws.Range("N10").Value = 1.736302
$ws.Range("O10").Value = 0.006489479200764093
$ws.Range("P10").Value = 0.004335698275101748
$ws.Range("Q10").Value = 88.85699983350999
$ws.Range("R10").Value = 533.14199900106
$ws.Range("S10").Value = 0.0037755300626405
$ws.Range("T10").Value = 0.002531101114312043
$ws.Range("G11").Value = 102.35201
$ws.Range("H11").Value = 307.05603
$ws.Range("I11").Value = 0.5817924591230612
$ws.Range("J11").Value = 0.583781654929077
$ws.Range("M11").Value = 48.13513433333333
$ws.Range("N11").Value = 144.405403
$ws.Range("O11").Value = 0.3598129277996019
$ws.Range("P11").Value = 0.3605929479448119
$ws.Range("Q11").Value = 4926.727750636675
$ws.Range("R11").Value = 44340.54975573008
$ws.Range("S11").Value = 0.2093364480887989
$ws.Range("T11").Value = 0.2105075479069768
$ws.Range("G12").Value = 1.79836
$ws.Range("H12").Value = 3.59672
$ws.Range("I12").Value = 0.01022229350247785
$ws.Range("J12").Value = 0.006838162904394061
$ws.Range("M12").Value = 0.346153
$ws.Range("N12").Value = 1.038459
$ws.Range("O12").Value = 0.002587513800919532
$ws.Range("P12").Value = 0.002593123140481257
$ws.Range("Q12").Value = 0.62250770908
$ws.Range("R12").Value = 3.73504625448
$ws.Range("S12").Value = 0.00002645032551471151
$ws.Range("T12").Value = 0.00001773219846576476
$ws.Range("G13").Value = 1.79836
$ws.Range("H13").Value = 3.59672
$ws.Range("I13").Value = 0.01022229350247785
$ws.Range("J13").Value = 0.006838162904394061
$ws.Range("O13").Value = 0.0001556606107424992
$ws.Range("P13").Value = 0.000155998059463248
$ws.Range("Q13").Value = 0.03744904864
$ws.Range("R13").Value = 0.22469429184
$ws.Range("S13").Value = 0.000001591208449784784
$ws.Range("T13").Value = 0.000001066740143379041
$ws.Range("G14").Value = 1.79836
$ws.Range("H14").Value = 3.59672
$ws.Range("I14").Value = 0.01022229350247785
$ws.Range("J14").Value = 0.006838162904394061
$ws.Range("M14").Value = 84.40796133333333
$ws.Range("N14").Value = 253.223884
$ws.Range("O14").Value = 0.630954418587972
$ws.Range("P14").Value = 0.6323222325801418
$ws.Range("Q14").Value = 151.7959013434133
$ws.Range("R14").Value = 910.77540806048
$ws.Range("S14").Value = 0.006449801253491518
$ws.Range("T14").Value = 0.004323922434453159
$ws.Range("G15").Value = 1.79836
$ws.Range("H15").Value = 3.59672
$ws.Range("I15").Value = 0.01022229350247785
$ws.Range("J15").Value = 0.006838162904394061
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.868151
$ws.Range("N15").Value = 1.736302
$ws.Range("O15").Value = 0.006489479200764093
$ws.Range("P15").Value = 0.004335698275101748
$ws.Range("Q15").Value = 1.56124803236
$ws.Range("R15").Value = 6.24499212944
$ws.Range("S15").Value = 0.00006633736106843597
$ws.Range("T15").Value = 0.00002964821110944609
$ws.Range("G16").Value = 1.79836
$ws.Range("H16").Value = 3.59672
$ws.Range("I16").Value = 0.01022229350247785
$ws.Range("J16").Value = 0.006838162904394061
$ws.Range("M16").Value = 48.13513433333333
$ws.Range("N16").Value = 144.405403
$ws.Range("O16").Value = 0.3598129277996019
$ws.Range("P16").Value = 0.3605929479448119
$ws.Range("Q16").Value = 86.56430017969332
$ws.Range("R16").Value = 519.3858010781599
$ws.Range("S16").Value = 0.003678113353953404
$ws.Range("T16").Value = 0.002465793320222311
$ws.Range("G17").Value = 71.02944933333333
$ws.Range("H17").Value = 213.088348
$ws.Range("I17").Value = 0.4037477915460271
$ws.Range("J17").Value = 0.4051282381314676
$ws.Range("M17").Value = 0.346153
$ws.Range("N17").Value = 1.038459
$ws.Range("O17").Value = 0.002587513800919532
$ws.Range("P17").Value = 0.002593123140481257
$ws.Range("Q17").Value = 24.58705697508133
$ws.Range("R17").Value = 221.283512775732
$ws.Range("S17").Value = 0.001044702982716128
$ws.Range("T17").Value = 0.00105054740916111
$ws.Range("G18").Value = 71.02944933333333
$ws.Range("H18").Value = 213.088348
$ws.Range("I18").Value = 0.4037477915460271
$ws.Range("J18").Value = 0.4051282381314676
$ws.Range("O18").Value = 0.0001556606107424992
$ws.Range("P18").Value = 0.000155998059463248
$ws.Range("Q18").Value = 1.479117252917333
$ws.Range("R18").Value = 13.312055276256
$ws.Range("S18").Value = 0.00006284762781798985
$ws.Range("T18").Value = 0.00006319921898227359
$ws.Range("G19").Value = 71.02944933333333
$ws.Range("H19").Value = 213.088348
$ws.Range("I19").Value = 0.4037477915460271
$ws.Range("J19").Value = 0.4051282381314676
$ws.Range("M19").Value = 84.40796133333333
$ws.Range("N19").Value = 253.223884
$ws.Range("O19").Value = 0.630954418587972
$ws.Range("P19").Value = 0.6323222325801418
$ws.Range("Q19").Value = 5995.451012855959
$ws.Range("R19").Value = 53959.05911570363
$ws.Range("S19").Value = 0.2547464530711013
$ws.Range("T19").Value = 0.2561715920165489
$ws.Range("G20").Value = 71.02944933333333
$ws.Range("H20").Value = 213.088348
$ws.Range("I20").Value = 0.4037477915460271
$ws.Range("J20").Value = 0.4051282381314676
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 0.868151
$ws.Range("N20").Value = 1.736302
$ws.Range("O20").Value = 0.006489479200764093
$ws.Range("P20").Value = 0.004335698275101748
$ws.Range("Q20").Value = 61.66428746818266
$ws.Range("R20").Value = 369.985724809096
$ws.Range("S20").Value = 0.00262011289559238
$ws.Range("T20").Value = 0.001756513803261615
$ws.Range("G21").Value = 71.02944933333333
$ws.Range("H21").Value = 213.088348
$ws.Range("I21").Value = 0.4037477915460271
$ws.Range("J21").Value = 0.4051282381314676
$ws.Range("M21").Value = 48.13513433333333
$ws.Range("N21").Value = 144.405403
$ws.Range("O21").Value = 0.3598129277996019
$ws.Range("P21").Value = 0.3605929479448119
$ws.Range("Q21").Value = 3419.012085282693
$ws.Range("R21").Value = 30771.10876754424
$ws.Range("S21").Value = 0.1452736749687994
$ws.Range("T21").Value = 0.1460863856835137
